$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 21.10438315076156
$ws.Range("C2").Value = 3.322157870623979
$ws.Range("D2").Value = 3.815026433025047
$ws.Range("F2").Value = 68.51717879731122
$ws.Range("G2").Value = 3.827825372342945
$ws.Range("J2").Value = 11.52299371725909
$ws.Range("K2").Value = 17.16141890428164
$ws.Range("L2").Value = 11.3718091587387

$ws.Range("B3").Value = 21.11827868237029
$ws.Range("C3").Value = 3.182415127181045
$ws.Range("D3").Value = 3.77283583236601
$ws.Range("F3").Value = 67.67545219562973
$ws.Range("G3").Value = 3.831506445824125
$ws.Range("J3").Value = 11.51673320975796
$ws.Range("K3").Value = 17.19030646397591
$ws.Range("L3").Value = 11.41862372721977

$ws.Range("B4").Value = 21.13358792487648
$ws.Range("C4").Value = 3.092480779927253
$ws.Range("D4").Value = 3.748175983357865
$ws.Range("F4").Value = 67.15613894465271
$ws.Range("G4").Value = 3.833882159496253
$ws.Range("J4").Value = 11.51343679520299
$ws.Range("K4").Value = 17.21340605080329
$ws.Range("L4").Value = 11.44949937710048

$ws.Range("B5").Value = 21.14152853120866
$ws.Range("C5").Value = 3.054809695380171
$ws.Range("D5").Value = 3.738449963671115
$ws.Range("F5").Value = 66.94403437763405
$ws.Range("G5").Value = 3.834879443448906
$ws.Range("J5").Value = 11.51223138172048
$ws.Range("K5").Value = 17.22416451996093
$ws.Range("L5").Value = 11.46261781923382

$ws.Range("B6").Value = 21.14294977005848
$ws.Range("C6").Value = 3.048493255434859
$ws.Range("D6").Value = 3.736854750424613
$ws.Range("F6").Value = 66.90878982881756
$ws.Range("G6").Value = 3.835046806279481
$ws.Range("J6").Value = 11.51203955845483
$ws.Range("K6").Value = 17.22603209535286
$ws.Range("L6").Value = 11.46482853443361

$ws.Range("B7").Value = 21.13368812734735
$ws.Range("C7").Value = 3.091976848617787
$ws.Range("D7").Value = 3.74804349430644
$ws.Range("F7").Value = 67.15328018037408
$ws.Range("G7").Value = 3.833895490994919
$ws.Range("J7").Value = 11.51341997993826
$ws.Range("K7").Value = 17.21354570143606
$ws.Range("L7").Value = 11.44967412472681

$ws.Range("B8").Value = 21.10776666965025
$ws.Range("C8").Value = 3.27484137273966
$ws.Range("D8").Value = 3.800226855091859
$ws.Range("F8").Value = 68.22755569055819
$ws.Range("G8").Value = 3.829070698537482
$ws.Range("J8").Value = 11.52072157933398
$ws.Range("K8").Value = 17.17026457670861
$ws.Range("L8").Value = 11.38750872223888

$ws.Range("B9").Value = 21.11076217555098
$ws.Range("C9").Value = 3.600050706767215
$ws.Range("D9").Value = 3.911963204974226
$ws.Range("F9").Value = 70.30807762262626
$ws.Range("G9").Value = 3.820520733592817
$ws.Range("J9").Value = 11.5393758825632
$ws.Range("K9").Value = 17.12805805455895
$ws.Range("L9").Value = 11.28249488535362

$ws.Range("B10").Value = 21.1457782351442
$ws.Range("C10").Value = 3.817987084996174
$ws.Range("D10").Value = 3.999102783868744
$ws.Range("F10").Value = 71.8122709487286
$ws.Range("G10").Value = 3.814787485306405
$ws.Range("J10").Value = 11.55571359603225
$ws.Range("K10").Value = 17.12318249928035
$ws.Range("L10").Value = 11.21561422746843

$ws.Range("B11").Value = 21.16880881902725
$ws.Range("C11").Value = 3.912477713400484
$ws.Range("D11").Value = 4.039683850802437
$ws.Range("F11").Value = 72.48952451255221
$ws.Range("G11").Value = 3.812296819165616
$ws.Range("J11").Value = 11.56371494097251
$ws.Range("K11").Value = 17.12665093280735
$ws.Range("L11").Value = 11.18741420755264

$ws.Range("B12").Value = 21.17854696324843
$ws.Range("C12").Value = 3.947584241158272
$ws.Range("D12").Value = 4.055172516015352
$ws.Range("F12").Value = 72.74483084958548
$ws.Range("G12").Value = 3.811370434348276
$ws.Range("J12").Value = 11.56682637344614
$ws.Range("K12").Value = 17.12878174084353
$ws.Range("L12").Value = 11.17705512057089

$ws.Range("B13").Value = 21.17640453068202
$ws.Range("C13").Value = 3.940053531244382
$ws.Range("D13").Value = 4.051831572519639
$ws.Range("F13").Value = 72.68989972040069
$ws.Range("G13").Value = 3.81156920347534
$ws.Range("J13").Value = 11.56615265220604
$ws.Range("K13").Value = 17.128286492527
$ws.Range("L13").Value = 11.17927192196016

$ws.Range("B14").Value = 21.16958962927862
$ws.Range("C14").Value = 3.915379508727648
$ws.Range("D14").Value = 4.040955757214239
$ws.Range("F14").Value = 72.51055269078141
$ws.Range("G14").Value = 3.812220269322891
$ws.Range("J14").Value = 11.56396928922781
$ws.Range("K14").Value = 17.1268098584272
$ws.Range("L14").Value = 11.18655555412738

$ws.Range("B15").Value = 21.16554759683701
$ws.Range("C15").Value = 3.900177883601669
$ws.Range("D15").Value = 4.034309416527988
$ws.Range("F15").Value = 72.40054276787242
$ws.Range("G15").Value = 3.812621247534687
$ws.Range("J15").Value = 11.56264251854006
$ws.Range("K15").Value = 17.12601180350141
$ws.Range("L15").Value = 11.19105861095114

$ws.Range("B16").Value = 21.14441570064728
$ws.Range("C16").Value = 3.811717756558672
$ws.Range("D16").Value = 3.996468472507136
$ws.Range("F16").Value = 71.76785712659473
$ws.Range("G16").Value = 3.814952609685427
$ws.Range("J16").Value = 11.55520211377847
$ws.Range("K16").Value = 17.12307024424137
$ws.Range("L16").Value = 11.21750190616477

$ws.Range("B17").Value = 21.13326823538261
$ws.Range("C17").Value = 3.756254419070119
$ws.Range("D17").Value = 3.973485494948566
$ws.Range("F17").Value = 71.37782367719616
$ws.Range("G17").Value = 3.816412821216276
$ws.Range("J17").Value = 11.5507832150787
$ws.Range("K17").Value = 17.12272212719691
$ws.Range("L17").Value = 11.23429356221573

$ws.Range("B18").Value = 21.12752527768741
$ws.Range("C18").Value = 3.723916005471954
$ws.Range("D18").Value = 3.960355691287095
$ws.Range("F18").Value = 71.15283535275644
$ws.Range("G18").Value = 3.81726375457576
$ws.Range("J18").Value = 11.54829518392065
$ws.Range("K18").Value = 17.12305719140376
$ws.Range("L18").Value = 11.24416105421864

$ws.Range("B19").Value = 21.12569578144353
$ws.Range("C19").Value = 3.712891874168975
$ws.Range("D19").Value = 3.955925936732247
$ws.Range("F19").Value = 71.07655082104314
$ws.Range("G19").Value = 3.817553768561951
$ws.Range("J19").Value = 11.54746199609015
$ws.Range("K19").Value = 17.12326257207567
$ws.Range("L19").Value = 11.24753798446513

$ws.Range("B20").Value = 21.13438571394651
$ws.Range("C20").Value = 3.762203902750159
$ws.Range("D20").Value = 3.975922919242834
$ws.Range("F20").Value = 71.41941189534306
$ws.Range("G20").Value = 3.816256235502274
$ws.Range("J20").Value = 11.55124806923232
$ws.Range("K20").Value = 17.1227037841243
$ws.Range("L20").Value = 11.23248439380568

$ws.Range("B21").Value = 21.17156377089657
$ws.Range("C21").Value = 3.922645231264877
$ws.Range("D21").Value = 4.044147061142806
$ws.Range("F21").Value = 72.56326375409918
$ws.Range("G21").Value = 3.812028581054344
$ws.Range("J21").Value = 11.56460838705535
$ws.Range("K21").Value = 17.12722140414986
$ws.Range("L21").Value = 11.18440750119885

$ws.Range("B22").Value = 21.20178657266246
$ws.Range("C22").Value = 4.023566343403381
$ws.Range("D22").Value = 4.089437171369379
$ws.Range("F22").Value = 73.30404023846522
$ws.Range("G22").Value = 3.809363300267107
$ws.Range("J22").Value = 11.57381496726322
$ws.Range("K22").Value = 17.13493761446043
$ws.Range("L22").Value = 11.1548495446088

$ws.Range("B23").Value = 21.18511569623125
$ws.Range("C23").Value = 3.970064805404743
$ws.Range("D23").Value = 4.065205373939103
$ws.Range("F23").Value = 72.90934324369988
$ws.Range("G23").Value = 3.810776903559242
$ws.Range("J23").Value = 11.56885791975668
$ws.Range("K23").Value = 17.13038374590698
$ws.Range("L23").Value = 11.17045477518242

$ws.Range("B24").Value = 21.13387842745466
$ws.Range("C24").Value = 3.759515547996048
$ws.Range("D24").Value = 3.974820698865317
$ws.Range("F24").Value = 71.40061219576475
$ws.Range("G24").Value = 3.816326992265362
$ws.Range("J24").Value = 11.55103774527275
$ws.Range("K24").Value = 17.12271040989922
$ws.Range("L24").Value = 11.23330165302668

$ws.Range("B25").Value = 21.10418441855446
$ws.Range("C25").Value = 3.515721559698256
$ws.Range("D25").Value = 3.880797301617488
$ws.Range("F25").Value = 69.74891591301208
$ws.Range("G25").Value = 3.82273689839097
$ws.Range("J25").Value = 11.53386622927416
$ws.Range("K25").Value = 17.13489281124436
$ws.Range("L25").Value = 11.30909811784316
